$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154 (shifts the existing rows 154..245 down to 155..246)
$ws.Rows.Item(154).Insert()

# Populate the new row 154 with the new weekly price record
$ws.Cells.Item(154, 1).Value  = 4
$ws.Cells.Item(154, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(154, 3).Value  = "Los Lagos"
$ws.Cells.Item(154, 4).Value  = 44452
$ws.Cells.Item(154, 5).Value  = 10
$ws.Cells.Item(154, 6).Value  = 100112006
$ws.Cells.Item(154, 7).Value  = "Repollo"
$ws.Cells.Item(154, 8).Value  = "Crespo record"
$ws.Cells.Item(154, 9).Value  = "Segunda"
$ws.Cells.Item(154, 10).Value = 500
$ws.Cells.Item(154, 11).Value = 1000
$ws.Cells.Item(154, 12).Value = 1000
$ws.Cells.Item(154, 13).Value = 1000
$ws.Cells.Item(154, 14).Value = "$/unidad"
$ws.Cells.Item(154, 15).Value = "Región del Maule"
$ws.Cells.Item(154, 16).Value = 1000
$ws.Cells.Item(154, 17).Value = 1
$ws.Cells.Item(154, 18).Value = "Hortaliza"
